$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update patient_id values (shared string replacements in the diff)
$ws.Range("A81").Value = "UCSF-PDGM-0298"
$ws.Range("A82").Value = "UCSF-PDGM-0307"

$ws.Range("A131").Value = "UCSF-PDGM-0540"
$ws.Range("A132").Value = "UCSF-PDGM-0281"

$ws.Range("A137").Value = "UCSF-PDGM-0269"

# Apply a new font-based style to rows 81,82 / 131,132 / 137 (A:B)
$ws.Range("A81:B82").Font.Size = 11
$ws.Range("A131:B132").Font.Size = 11
$ws.Range("A137:B137").Font.Size = 11

# Restore the selection/view as recorded in the sheet
$ws.Application.ActiveWindow.ScrollRow = 80
$ws.Range("A82:B82").Select
